$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1): make bold (adds a new bold font + cell style, keeps
#    the existing centered alignment).
# ---------------------------------------------------------------------------
$ws.Range("A1:R1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) New data rows 5-7 (additional "W-" boson measurements), mirroring the
#    layout of existing rows 2-4.
# ---------------------------------------------------------------------------

# Make sure the new rows carry the same centered-alignment style ("s=1") as
# the rest of the data cells (columns A-P only; Q and R stay unstyled, same
# as in the existing rows 2-4), even for cells that end up blank.
$ws.Range("A5:P7").HorizontalAlignment = -4108

# Row 5
$ws.Range("A5").Value = "STAR"
$ws.Range("B5").Value = "pp"
$ws.Range("C5").Value = 510
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 0.27
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = "W-"
$ws.Range("J5").Value = -0.012
$ws.Range("K5").Value = 0.101
$ws.Range("L5").Value = 0.019
$ws.Range("M5").Value = "A_LL"
$ws.Range("N5").Value = "eta"
$ws.Range("O5").Value = 0.005

# Row 6
$ws.Range("A6").Value = "STAR"
$ws.Range("B6").Value = "pp"
$ws.Range("C6").Value = 510
$ws.Range("D6").Value = 0.5
$ws.Range("E6").Value = 1.1
$ws.Range("F6").Value = 0.74
$ws.Range("G6").Value = 25
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = "W-"
$ws.Range("J6").Value = -0.028
$ws.Range("K6").Value = 0.092
$ws.Range("L6").Value = 0.02
$ws.Range("M6").Value = "A_LL"
$ws.Range("N6").Value = "eta"
$ws.Range("O6").Value = 0.005

# Row 7
$ws.Range("A7").Value = "STAR"
$ws.Range("B7").Value = "pp"
$ws.Range("C7").Value = 510
$ws.Range("D7").Value = 1.1
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1.27
$ws.Range("G7").Value = 25
$ws.Range("H7").Value = 50
$ws.Range("I7").Value = "W-"
$ws.Range("J7").Value = -0.147
$ws.Range("K7").Value = 0.26
$ws.Range("L7").Value = 0.038
$ws.Range("M7").Value = "A_LL"
$ws.Range("N7").Value = "eta"
$ws.Range("O7").Value = 0.005

foreach ($r in 5..7) {
    $ws.Range("P$r").Formula = "=0.033*J$r"
    $ws.Range("Q$r").Formula = "=SQRT(L$r*L$r-O$r*O$r-P$r*P$r)"
    $ws.Range("R$r").Formula = "=SQRT(O$r*O$r+P$r*P$r)"
}

# ---------------------------------------------------------------------------
# 3) Selection moves to P4 (matches the saved workbook view state).
# ---------------------------------------------------------------------------
$ws.Range("P4").Select()
